$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.433935198812251
$ws.Range("C2").Value = 0.4883748341653806
$ws.Range("D2").Value = 0.2913197292499859
$ws.Range("E2").Value = 0.5397404276594314
$ws.Range("F2").Value = 0.3340733006221409
$ws.Range("G2").Value = 13

$ws.Range("B3").Value = 0.3043483804216044
$ws.Range("C3").Value = 0.3187202035556604
$ws.Range("D3").Value = 0.1391823142959026
$ws.Range("E3").Value = 0.3730714600393638
$ws.Range("F3").Value = 0.2253588111853846
$ws.Range("G3").Value = 12

$ws.Range("B4").Value = 0.2239330230754586
$ws.Range("C4").Value = 0.2420726978122317
$ws.Range("D4").Value = 0.08325831886165148
$ws.Range("E4").Value = 0.2885451764657512
$ws.Range("F4").Value = 0.190849553422929
$ws.Range("G4").Value = 11

$ws.Range("B5").Value = 0.3439814776516536
$ws.Range("C5").Value = 0.3439814776516536
$ws.Range("D5").Value = 0.1481954533256283
$ws.Range("E5").Value = 0.3849616257831789
$ws.Range("F5").Value = 0.1833199959169479
$ws.Range("G5").Value = 9

$ws.Range("B6").Value = 0.3607780907168635
$ws.Range("C6").Value = 0.3638782949054072
$ws.Range("D6").Value = 0.166400635697807
$ws.Range("E6").Value = 0.40792234027791
$ws.Range("F6").Value = 0.2006650192526551
$ws.Range("G6").Value = 10

$ws.Range("B7").Value = 0.3167374129797205
$ws.Range("C7").Value = 0.3281249448458587
$ws.Range("D7").Value = 0.1360088531859135
$ws.Range("E7").Value = 0.3687937813818361
$ws.Range("F7").Value = 0.2003672813995112
$ws.Range("G7").Value = 9

$ws.Range("B8").Value = 0.3332309007927744
$ws.Range("C8").Value = 0.3363297904028241
$ws.Range("D8").Value = 0.1472445396808239
$ws.Range("E8").Value = 0.3837245622589515
$ws.Range("F8").Value = 0.2034044709093121
$ws.Range("G8").Value = 8

$ws.Range("B9").Value = 0.2929253822463562
$ws.Range("C9").Value = 0.295209788697302
$ws.Range("D9").Value = 0.1182577732422086
$ws.Range("E9").Value = 0.3438862795201469
$ws.Range("F9").Value = 0.1973397892307624
$ws.Range("G9").Value = 6

$ws.Range("B10").Value = 0.3047033996258702
$ws.Range("C10").Value = 0.3047033996258702
$ws.Range("D10").Value = 0.1225151295288574
$ws.Range("E10").Value = 0.35002161294534
$ws.Range("F10").Value = 0.1886932996753026
$ws.Range("G10").Value = 6

$ws.Range("B11").Value = 0.2886753471776797
$ws.Range("C11").Value = 0.2886753471776797
$ws.Range("D11").Value = 0.1109408868566417
$ws.Range("E11").Value = 0.3330778990816439
$ws.Range("F11").Value = 0.1857667582900928
$ws.Range("G11").Value = 5
